$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column F (old F shifts to G) ---
$ws.Columns.Item(6).Insert()

# --- Fill in new cell values, in the order the author likely typed them so
#     that shared-string indices come out in the same order as the target file ---
$ws.Range("F3").Value = "pagGeldVerdienen"
$ws.Range("C2").Value = '//android.view.View[@content-desc=" Profil"]/android.widget.TextView[@text="Profil"]'
$ws.Range("E2").Value = '//android.view.View[@content-desc=" Benachrichtigungen"]/android.widget.TextView[@text="Benachrichtigungen"]'
$ws.Range("D2").Value = '//android.view.View[@content-desc=" Abwesenheiten"]/android.widget.TextView[@text="Abwesenheiten"]'
$ws.Range("F2").Value = '//android.view.View[@content-desc=" Geld verdienen"]/android.widget.TextView[@text="Geld verdienen"]'

# --- Give the (still empty) F1 header cell an orange fill (new style) ---
$ws.Range("F1").Interior.Color = 49407

# --- Resize the new/expanded columns to their best-fit-like widths ---
$pad = 0.8333333333333334
$ws.Columns.Item(3).ColumnWidth = 83.44140625 - $pad
$ws.Columns.Item(4).ColumnWidth = 92 - $pad
$ws.Columns.Item(5).ColumnWidth = 99.109375 - $pad
$ws.Columns.Item(6).ColumnWidth = 91.33203125 - $pad

# --- Nudge the picture so Excel recalculates its bottom-right cell anchor
#     to reflect the new (wider) columns, while leaving its top-left anchor
#     untouched ---
$shp = $ws.Shapes.Item(1)
$shp.Width = 1159.2273

# --- Update the view: scroll so column E is left-most visible, and select F2 ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 5
$win.ScrollRow = 1
$ws.Range("F2").Select() | Out-Null
